$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '27.809.14'
$ws.Range('E2').Value = '  +1.43%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.863.42'
$ws.Range('E3').Value = '  +1.44%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.033'
$ws.Range('E4').Value = '  +0.51%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '322.42'
$ws.Range('E5').Value = '  +1.33%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '1.029'
$ws.Range('E6').Value = '  +0.37%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4420'
$ws.Range('E7').Value = '  +1.38%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3810'
$ws.Range('E8').Value = '  +2.30%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.07437'
$ws.Range('E9').Value = '  +1.37%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.8867'
$ws.Range('E10').Value = '  +1.61%  '
$ws.Range('E11').Value = '  +1.87%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.879.49'
$ws.Range('E12').Value = '  -4.80%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '5.551'
$ws.Range('E13').Value = '  +1.51%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.771'
$ws.Range('E14').Value = '  +1.35%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.07186'
$ws.Range('E15').Value = '  +0.53%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '84.40'
$ws.Range('E16').Value = '  +2.65%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '1.035'
$ws.Range('E17').Value = '  +0.43%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.000009119'
$ws.Range('E18').Value = '  +1.52%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '1.030'
$ws.Range('E19').Value = '  +0.52%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '15.52'
$ws.Range('E20').Value = '  +0.75%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '27.803.13'
$ws.Range('E21').Value = '  +1.33%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.314'
$ws.Range('E22').Value = '  +1.23%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '11.31'
$ws.Range('E23').Value = '  +1.10%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.092.67'
$ws.Range('E24').Value = '  -3.30%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.022'
$ws.Range('E25').Value = '  +6.86%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '158.23'
$ws.Range('E26').Value = '  +0.79%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '18.90'
$ws.Range('E27').Value = '  +1.95%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '5.384'
$ws.Range('E28').Value = '  +2.59%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.984'
$ws.Range('E29').Value = '  +3.06%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '118.84'
$ws.Range('E30').Value = '  +3.07%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.09031'
$ws.Range('E31').Value = '  -0.18%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.236'
$ws.Range('E32').Value = '  +3.27%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.7782'
$ws.Range('E33').Value = '  +2.57%  '
$ws.Range('B34').Value = 'Filecoin'
$ws.Range('C34').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '4.598'
$ws.Range('E34').Value = '  +3.05%  '
$ws.Range('B35').Value = 'HuobiToken'
$ws.Range('C35').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.994'
$ws.Range('E35').Value = '  +4.65%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.032'
$ws.Range('E36').Value = '  +0.41%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.143'
$ws.Range('E37').Value = '  -0.68%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.01988'
$ws.Range('E38').Value = '  +1.62%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.05360'
$ws.Range('E39').Value = '  +2.13%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.881'
$ws.Range('E40').Value = '  +3.42%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.5211'
$ws.Range('E41').Value = '  +0.79%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.1694'
$ws.Range('E42').Value = '  +2.03%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '6.903'
$ws.Range('E43').Value = '  +5.97%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '8.732'
$ws.Range('E44').Value = '  +3.07%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '111.44'
$ws.Range('E45').Value = '  +2.52%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '10.72'
$ws.Range('E46').Value = '  +2.29%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.06660'
$ws.Range('E47').Value = '  +5.87%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.031'
$ws.Range('E48').Value = '  +0.30%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.713'
$ws.Range('E49').Value = '  +2.69%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.4730'
$ws.Range('E50').Value = '  +2.21%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.925'
$ws.Range('E51').Value = '  +2.97%  '
